$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I7').Value = 'sv'
$ws.Range('J7').Value = 'Statement-opinion'
$ws.Range('I10').Value = 'ba'
$ws.Range('J10').Value = 'Appreciation'
$ws.Range('I17').Value = 'aa'
$ws.Range('J17').Value = 'Agree/Accept'
$ws.Range('I19').Value = 'sd'
$ws.Range('J19').Value = 'Statement-non-opinion'
$ws.Range('I45').Value = 'aa'
$ws.Range('J45').Value = 'Agree/Accept'
$ws.Range('I47').Value = 'aa'
$ws.Range('J47').Value = 'Agree/Accept'
$ws.Range('I48').Value = 'aa'
$ws.Range('J48').Value = 'Agree/Accept'
$ws.Range('I52').Value = 'aa'
$ws.Range('J52').Value = 'Agree/Accept'
$ws.Range('I56').Value = 'b'
$ws.Range('J56').Value = 'Acknowledge (Backchannel)'
$ws.Range('I72').Value = 'b'
$ws.Range('J72').Value = 'Acknowledge (Backchannel)'
$ws.Range('I74').Value = 'sv'
$ws.Range('J74').Value = 'Statement-opinion'
$ws.Range('I75').Value = 'sd'
$ws.Range('J75').Value = 'Statement-non-opinion'
$ws.Range('I91').Value = '%'
$ws.Range('J91').Value = 'Uninterpretable'
$ws.Range('I97').Value = 'sd'
$ws.Range('J97').Value = 'Statement-non-opinion'
$ws.Range('I100').Value = 'sd'
$ws.Range('J100').Value = 'Statement-non-opinion'
$ws.Range('I107').Value = 'sv'
$ws.Range('J107').Value = 'Statement-opinion'
$ws.Range('I110').Value = 'aa'
$ws.Range('J110').Value = 'Agree/Accept'
$ws.Range('I116').Value = 'sd'
$ws.Range('J116').Value = 'Statement-non-opinion'
$ws.Range('I126').Value = 'b'
$ws.Range('J126').Value = 'Acknowledge (Backchannel)'
$ws.Range('I127').Value = 'b'
$ws.Range('J127').Value = 'Acknowledge (Backchannel)'
$ws.Range('I149').Value = 'b'
$ws.Range('J149').Value = 'Acknowledge (Backchannel)'
$ws.Range('I156').Value = 'ba'
$ws.Range('J156').Value = 'Appreciation'
$ws.Range('I164').Value = 'b'
$ws.Range('J164').Value = 'Acknowledge (Backchannel)'
$ws.Range('I167').Value = 'sd'
$ws.Range('J167').Value = 'Statement-non-opinion'
$ws.Range('I168').Value = 'sd'
$ws.Range('J168').Value = 'Statement-non-opinion'
$ws.Range('I179').Value = 'aa'
$ws.Range('J179').Value = 'Agree/Accept'
$ws.Range('I186').Value = 'sv'
$ws.Range('J186').Value = 'Statement-opinion'
$ws.Range('I188').Value = 'aa'
$ws.Range('J188').Value = 'Agree/Accept'
$ws.Range('I196').Value = 'qy'
$ws.Range('J196').Value = 'Yes-No-Question'
$ws.Range('I202').Value = 'sv'
$ws.Range('J202').Value = 'Statement-opinion'
$ws.Range('I213').Value = 'sv'
$ws.Range('J213').Value = 'Statement-opinion'
$ws.Range('I216').Value = 'sv'
$ws.Range('J216').Value = 'Statement-opinion'
$ws.Range('I223').Value = 'sd'
$ws.Range('J223').Value = 'Statement-non-opinion'
$ws.Range('I224').Value = 'sd'
$ws.Range('J224').Value = 'Statement-non-opinion'
$ws.Range('I250').Value = 'sd'
$ws.Range('J250').Value = 'Statement-non-opinion'
$ws.Range('I252').Value = 'aa'
$ws.Range('J252').Value = 'Agree/Accept'
$ws.Range('I266').Value = 'sd'
$ws.Range('J266').Value = 'Statement-non-opinion'
$ws.Range('I269').Value = 'ba'
$ws.Range('J269').Value = 'Appreciation'
$ws.Range('I280').Value = 'b'
$ws.Range('J280').Value = 'Acknowledge (Backchannel)'
$ws.Range('I284').Value = 'ba'
$ws.Range('J284').Value = 'Appreciation'
$ws.Range('I288').Value = 'sd'
$ws.Range('J288').Value = 'Statement-non-opinion'
$ws.Range('I290').Value = 'sv'
$ws.Range('J290').Value = 'Statement-opinion'
$ws.Range('I294').Value = 'aa'
$ws.Range('J294').Value = 'Agree/Accept'
$ws.Range('I307').Value = 'sv'
$ws.Range('J307').Value = 'Statement-opinion'
$ws.Range('I310').Value = 'sd'
$ws.Range('J310').Value = 'Statement-non-opinion'
$ws.Range('I343').Value = 'aa'
$ws.Range('J343').Value = 'Agree/Accept'
$ws.Range('I349').Value = 'sd'
$ws.Range('J349').Value = 'Statement-non-opinion'
$ws.Range('I350').Value = 'sv'
$ws.Range('J350').Value = 'Statement-opinion'
$ws.Range('I352').Value = '%'
$ws.Range('J352').Value = 'Uninterpretable'
$ws.Range('I361').Value = 'sd'
$ws.Range('J361').Value = 'Statement-non-opinion'
$ws.Range('I370').Value = 'sd'
$ws.Range('J370').Value = 'Statement-non-opinion'
$ws.Range('I377').Value = 'sd'
$ws.Range('J377').Value = 'Statement-non-opinion'
$ws.Range('I379').Value = 'aa'
$ws.Range('J379').Value = 'Agree/Accept'
$ws.Range('I387').Value = 'ba'
$ws.Range('J387').Value = 'Appreciation'
$ws.Range('I404').Value = 'aa'
$ws.Range('J404').Value = 'Agree/Accept'
$ws.Range('I417').Value = '%'
$ws.Range('J417').Value = 'Uninterpretable'
$ws.Range('I418').Value = 'aa'
$ws.Range('J418').Value = 'Agree/Accept'
$ws.Range('I428').Value = '%'
$ws.Range('J428').Value = 'Uninterpretable'
$ws.Range('I429').Value = 'aa'
$ws.Range('J429').Value = 'Agree/Accept'
$ws.Range('I448').Value = 'sd'
$ws.Range('J448').Value = 'Statement-non-opinion'
$ws.Range('I451').Value = 'sv'
$ws.Range('J451').Value = 'Statement-opinion'
$ws.Range('I458').Value = 'sv'
$ws.Range('J458').Value = 'Statement-opinion'
$ws.Range('I471').Value = 'aa'
$ws.Range('J471').Value = 'Agree/Accept'
$ws.Range('I472').Value = 'sd'
$ws.Range('J472').Value = 'Statement-non-opinion'
$ws.Range('I481').Value = '%'
$ws.Range('J481').Value = 'Uninterpretable'
$ws.Range('I489').Value = 'sd'
$ws.Range('J489').Value = 'Statement-non-opinion'
$ws.Range('I491').Value = 'sd'
$ws.Range('J491').Value = 'Statement-non-opinion'
